# Re-sort the curvature calibration data (rows below the header) in
# ascending order of column A (the "time (s)" column), as if the data
# table had been re-sorted by timestamp after performing the needle
# calibration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (header is row 1, data starts row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # -4162 = xlUp

$dataRange = $ws.Range("A2:D" + $lastRow)
$keyRange = $ws.Range("A2:A" + $lastRow)

# 1 = xlAscending
$dataRange.Sort($keyRange, 1)
